$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block1 = New-Object 'object[,]' 24,7
$block1[0,0] = 0.8303843263742863
$block1[0,1] = 0.1067054068322477
$block1[0,2] = 0
$block1[0,3] = 0.1036435616348044
$block1[0,4] = 0.4443680307746263
$block1[0,5] = 1.430881998654868
$block1[0,6] = 1.332809327528906
$block1[1,0] = 0.7972251687634753
$block1[1,1] = 0.1050585938253974
$block1[1,2] = 0
$block1[1,3] = 0.1038382882349609
$block1[1,4] = 0.387822817061874
$block1[1,5] = 1.435355864311873
$block1[1,6] = 1.339798634065673
$block1[2,0] = 0.7772777148367993
$block1[2,1] = 0.1040301257229146
$block1[2,2] = 0
$block1[2,3] = 0.1039925699847881
$block1[2,4] = 0.3531389305168915
$block1[2,5] = 1.438836345955579
$block1[2,6] = 1.344599376255275
$block1[3,0] = 0.7692530707309402
$block1[3,1] = 0.1036066562302906
$block1[3,2] = 0
$block1[3,3] = 0.1040641896836974
$block1[3,4] = 0.3390132514313251
$block1[3,5] = 1.440439036433304
$block1[3,6] = 1.34668385030065
$block1[4,0] = 0.7679268843137663
$block1[4,1] = 0.1035360759253692
$block1[4,2] = 0
$block1[4,3] = 0.1040766110150404
$block1[4,4] = 0.336668177824194
$block1[4,5] = 1.44071629488117
$block1[4,6] = 1.34703771661033
$block1[5,0] = 0.7771690696268649
$block1[5,1] = 0.1040244323288064
$block1[5,2] = 0
$block1[5,3] = 0.1039935004256964
$block1[5,4] = 0.3529483938344953
$block1[5,5] = 1.438857214018881
$block1[5,6] = 1.344626969306475
$block1[6,0] = 0.8188656839572843
$block1[6,1] = 0.1061411763523452
$block1[6,2] = 0
$block1[6,3] = 0.1037035083077793
$block1[6,4] = 0.4248636149813336
$block1[6,5] = 1.432272296002878
$block1[6,6] = 1.33511358051345
$block1[7,0] = 0.9038922276211565
$block1[7,1] = 0.1101551040681272
$block1[7,2] = 0
$block1[7,3] = 0.103409567812804
$block1[7,4] = 0.5661985755041457
$block1[7,5] = 1.425184280816268
$block1[7,6] = 1.320496112909865
$block1[8,0] = 0.9683389370032671
$block1[8,1] = 0.1130214179392652
$block1[8,2] = 0
$block1[8,3] = 0.103360179571391
$block1[8,4] = 0.6702781546542269
$block1[8,5] = 1.423536407851046
$block1[8,6] = 1.312215575935511
$block1[9,0] = 0.9980853358762829
$block1[9,1] = 0.1143076019786378
$block1[9,2] = 0
$block1[9,3] = 0.1033736936665566
$block1[9,4] = 0.7176906081379002
$block1[9,5] = 1.42356170712597
$block1[9,6] = 1.308981936619816
$block1[10,0] = 1.009410957406118
$block1[10,1] = 0.1147921081971646
$block1[10,2] = 0
$block1[10,3] = 0.1033839683222109
$block1[10,4] = 0.7356546913071611
$block1[10,5] = 1.423682857504787
$block1[10,6] = 1.307834070970173
$block1[11,0] = 1.006969062290807
$block1[11,1] = 0.114687874186636
$block1[11,2] = 0
$block1[11,3] = 0.103381526393175
$block1[11,4] = 0.7317853510981394
$block1[11,5] = 1.423651801302682
$block1[11,6] = 1.308077876148204
$block1[12,0] = 0.9990158753309402
$block1[12,1] = 0.1143475135549394
$block1[12,2] = 0
$block1[12,3] = 0.1033744357078845
$block1[12,4] = 0.7191683204515869
$block1[12,5] = 1.423569437422657
$block1[12,6] = 1.3088859651577
$block1[13,0] = 0.9941522899871131
$block1[13,1] = 0.114138701683693
$block1[13,2] = 0
$block1[13,3] = 0.1033707635733911
$block1[13,4] = 0.7114413442032514
$block1[13,5] = 1.423533520965577
$block1[13,6] = 1.309390923288788
$block1[14,0] = 0.9664035417249011
$block1[14,1] = 0.1129370066142528
$block1[14,2] = 0
$block1[14,3] = 0.1033600189806911
$block1[14,4] = 0.6671810134426437
$block1[14,5] = 1.423550359627512
$block1[14,6] = 1.312437628755987
$block1[15,0] = 0.949490228473195
$block1[15,1] = 0.1121952707191625
$block1[15,2] = 0
$block1[15,3] = 0.1033626311069984
$block1[15,4] = 0.6400460337125793
$block1[15,5] = 1.42375926349942
$block1[15,6] = 1.314443230180615
$block1[16,0] = 0.9398025679213333
$block1[16,1] = 0.1117669760565008
$block1[16,2] = 0
$block1[16,3] = 0.1033675212450778
$block1[16,4] = 0.6244449056556647
$block1[16,5] = 1.423952350829339
$block1[16,6] = 1.315646989987812
$block1[17,0] = 0.9365294492675957
$block1[17,1] = 0.1116216762815228
$block1[17,2] = 0
$block1[17,3] = 0.1033697593195644
$block1[17,4] = 0.619163680173358
$block1[17,5] = 1.424030248762904
$block1[17,6] = 1.316063183737796
$block1[18,0] = 0.9512864982614815
$block1[18,1] = 0.1122744023534921
$block1[18,2] = 0
$block1[18,3] = 0.1033620025487778
$block1[18,4] = 0.642933953830422
$block1[18,5] = 1.423729476213879
$block1[18,6] = 1.314224536152594
$block1[19,0] = 1.00135025866868
$block1[19,1] = 0.1144475547303401
$block1[19,2] = 0
$block1[19,3] = 0.1033763785843451
$block1[19,4] = 0.7228739723491628
$block1[19,5] = 1.42359060061446
$block1[19,6] = 1.308646530002704
$block1[20,0] = 1.034426962160722
$block1[20,1] = 0.1158530102303246
$block1[20,2] = 0
$block1[20,3] = 0.103415826253805
$block1[20,4] = 0.7751780083420101
$block1[20,5] = 1.424150217511794
$block1[20,6] = 1.305447702524134
$block1[21,0] = 1.016740763836907
$block1[21,1] = 0.1151042474083397
$block1[21,2] = 0
$block1[21,3] = 0.1033920280779199
$block1[21,4] = 0.7472568307830727
$block1[21,5] = 1.42379198413812
$block1[21,6] = 1.307114112356146
$block1[22,0] = 0.9504742924330003
$block1[22,1] = 0.1122386327386593
$block1[22,2] = 0
$block1[22,3] = 0.1033622761645283
$block1[22,4] = 0.6416283278902171
$block1[22,5] = 1.423742715706268
$block1[22,6] = 1.314323249747815
$block1[23,0] = 0.8805423102667191
$block1[23,1] = 0.1090838144836397
$block1[23,2] = 0
$block1[23,3] = 0.1034597692701826
$block1[23,4] = 0.5279251897347166
$block1[23,5] = 1.426477285122616
$block1[23,6] = 1.324018508403185

$block2 = New-Object 'object[,]' 24,4
$block2[0,0] = 0.4386586759883073
$block2[0,1] = 0.1931854430375068
$block2[0,2] = 0.2099128288949856
$block2[0,3] = 2.581155204035216
$block2[1,0] = 0.407164507892162
$block2[1,1] = 0.190700560962263
$block2[1,2] = 0.203792270929231
$block2[1,3] = 2.604529354432923
$block2[2,0] = 0.3880170742764903
$block2[2,1] = 0.1892669514525096
$block2[2,2] = 0.2001396634215631
$block2[2,3] = 2.619638802877816
$block2[3,0] = 0.3802623926949309
$block2[3,1] = 0.18870595909695
$block2[3,2] = 0.1986777939036948
$block2[3,3] = 2.625986534335567
$block2[4,0] = 0.3789776433493017
$block2[4,1] = 0.1886142103971054
$block2[4,2] = 0.198436660501315
$block2[4,3] = 2.627052076531811
$block2[5,0] = 0.3879122970341911
$block2[5,1] = 0.1892592916477085
$block2[5,2] = 0.2001198403332154
$block2[5,3] = 2.61972363924956
$block2[6,0] = 0.4277601667052977
$block2[6,1] = 0.192309566081768
$block2[6,2] = 0.2077806243833571
$block2[6,3] = 2.589057366389191
$block2[7,0] = 0.5074036632164791
$block2[7,1] = 0.1990203353757352
$block2[7,2] = 0.223637336500154
$block2[7,3] = 2.534931568933732
$block2[8,0] = 0.5668310602073063
$block2[8,1] = 0.2043937332259844
$block2[8,2] = 0.2357935449480308
$block2[8,3] = 2.498827336642194
$block2[9,0] = 0.5940643584593772
$block2[9,1] = 0.2069341488516585
$block2[9,2] = 0.2414333025440882
$block2[9,3] = 2.483197281497603
$block2[10,0] = 0.6044054304394137
$block2[10,1] = 0.207909907817438
$block2[10,2] = 0.2435846644789734
$block2[10,3] = 2.477392769680993
$block2[11,0] = 0.6021770362104348
$block2[11,1] = 0.2076991495438989
$block2[11,2] = 0.2431206329436861
$block2[11,3] = 2.478637792456915
$block2[12,0] = 0.5949145567913945
$block2[12,1] = 0.2070141496717497
$block2[12,2] = 0.241609982065043
$block2[12,3] = 2.482717450571506
$block2[13,0] = 0.5904697705983324
$block2[13,1] = 0.2065963577691576
$block2[13,2] = 0.2406867078846417
$block2[13,3] = 2.485231240323913
$block2[14,0] = 0.5650552908875852
$block2[14,1] = 0.2042296391582283
$block2[14,2] = 0.2354271763502283
$block2[14,3] = 2.49986478064924
$block2[15,0] = 0.5495152283309892
$block2[15,1] = 0.2028022937218026
$block2[15,2] = 0.2322286946191596
$block2[15,3] = 2.509045389012549
$block2[16,0] = 0.5405957873341549
$block2[16,1] = 0.2019903641786698
$block2[16,2] = 0.2303993560941251
$block2[16,3] = 2.514400556412543
$block2[17,0] = 0.5375790536890861
$block2[17,1] = 0.2017170132576496
$block2[17,2] = 0.2297817526534871
$block2[17,3] = 2.516226555584108
$block2[18,0] = 0.551167551589856
$block2[18,1] = 0.2029533016527978
$block2[18,2] = 0.2325681085651681
$block2[18,3] = 2.508060363300082
$block2[19,0] = 0.5970469540850445
$block2[19,1] = 0.2072149777321215
$block2[19,2] = 0.2420532712201577
$block2[19,3] = 2.481516054675264
$block2[20,0] = 0.6271972825047385
$block2[20,1] = 0.2100804023751124
$block2[20,2] = 0.2483438923706061
$block2[20,3] = 2.464833797023307
$block2[21,0] = 0.6110904324286821
$block2[21,1] = 0.2085437520794926
$block2[21,2] = 0.2449781239716131
$block2[21,3] = 2.473676466071073
$block2[22,0] = 0.5504204902866832
$block2[22,1] = 0.2028850039568084
$block2[22,2] = 0.2324146298874936
$block2[22,3] = 2.508505453217339
$block2[23,0] = 0.4856974337952238
$block2[23,1] = 0.1971269890830385
$block2[23,2] = 0.2192586156575018
$block2[23,3] = 2.548930797048079

$ws.Range("B2:H25").Value2 = $block1
$ws.Range("K2:N25").Value2 = $block2
